$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A, rows 3 through 158 hold ALL-CAPS vendor codes like
# "PFIZER_JANSSEN" that should become Title_Case, e.g. "Pfizer_Janssen".
$lastRow = 158
for ($r = 3; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val.ToString().Split("_")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $p = $parts[$i]
            if ($p.Length -gt 0) {
                $parts[$i] = $p.Substring(0,1).ToUpper() + $p.Substring(1).ToLower()
            }
        }
        $newVal = [string]::Join("_", $parts)
        $cell.Value2 = $newVal
    }
}
